$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.588.25"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "1.922.45"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'247.15"
$ws.Range("E5").Value = "  +5.00%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4726"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").Value = "'0.2895"
$ws.Range("E8").Value = "  +4.88%  "
$ws.Range("D9").Value = "'0.06824"
$ws.Range("E9").Value = "  +6.71%  "
$ws.Range("D10").Value = "'105.41"
$ws.Range("E10").Value = "  +8.84%  "
$ws.Range("D11").Value = "'18.48"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "1.920.18"
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").Value = "'0.07693"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "'5.310"
$ws.Range("E14").Value = "  +7.02%  "
$ws.Range("D15").Value = "'0.6739"
$ws.Range("E15").Value = "  +8.50%  "
$ws.Range("D16").Value = "'290.01"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").Value = "30.597.69"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "'0.000007628"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'12.94"
$ws.Range("E20").Value = "  +3.17%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.541"
$ws.Range("E21").Value = "  +11.31%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.165.38"
$ws.Range("E22").Value = "  +4.42%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'6.335"
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("D26").Value = "'168.95"
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("D27").Value = "'21.24"
$ws.Range("E27").Value = "  +11.09%  "
$ws.Range("D28").Value = "'2.123"
$ws.Range("E28").Value = "  +10.37%  "
$ws.Range("D29").Value = "'0.1079"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "'1.395"
$ws.Range("E30").Value = "  +5.45%  "
$ws.Range("D31").Value = "'4.187"
$ws.Range("E31").Value = "  +5.47%  "
$ws.Range("D32").Value = "'4.127"
$ws.Range("E32").Value = "  +8.14%  "
$ws.Range("D33").Value = "'0.05073"
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").Value = "'0.7429"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("D35").Value = "'1.154"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").Value = "'0.02073"
$ws.Range("E36").Value = "  +9.78%  "
$ws.Range("D37").Value = "'2.750"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "'2.694"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "'2.068"
$ws.Range("E39").Value = "  +5.77%  "
$ws.Range("D40").Value = "'111.49"
$ws.Range("E40").Value = "  +4.94%  "
$ws.Range("D41").Value = "'0.8815"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").Value = "'0.4451"
$ws.Range("E42").Value = "  +10.67%  "
$ws.Range("D43").Value = "'5.893"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'67.13"
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").Value = "'7.263"
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("D47").Value = "'9.371"
$ws.Range("E47").Value = "  +5.44%  "
$ws.Range("D48").Value = "'48.62"
$ws.Range("E48").Value = "  +19.84%  "
$ws.Range("D49").Value = "'0.1235"
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").Value = "'0.4141"
$ws.Range("E50").Value = "  +12.09%  "
$ws.Range("D51").Value = "'34.98"
$ws.Range("E51").Value = "  +3.85%  "
